$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indicators")

$ws.Range("B2").Value  = "RAMSAR_A_1"
$ws.Range("B3").Value  = "RAMSAR_Ai_2"
$ws.Range("B4").Value  = "RAMSAR_Aii_3"
$ws.Range("B5").Value  = "RAMSAR_B_4"
$ws.Range("B6").Value  = "RAMSAR_Bi_5"
$ws.Range("B7").Value  = "RAMSAR_C_6"
$ws.Range("B8").Value  = "RAMSAR_Ci_7"
$ws.Range("B9").Value  = "RAMSAR_Cii_8"
$ws.Range("B10").Value = "RAMSAR_D_9"
$ws.Range("B11").Value = "RAMSAR_Di_10"
$ws.Range("B12").Value = "RAMSAR_E_11"
$ws.Range("B13").Value = "RAMSAR_Ei_12"
$ws.Range("B14").Value = "RAMSAR_F_13"
$ws.Range("B15").Value = "RAMSAR_Fi_14"
$ws.Range("B16").Value = "RAMSAR_G_15"
$ws.Range("B17").Value = "RAMSAR_Gi_16"
$ws.Range("B18").Value = "RAMSAR_Gii_17"
$ws.Range("B19").Value = "RAMSAR_H_18"
$ws.Range("B20").Value = "RAMSAR_Hi_19"

$ws.Columns.Item(2).ColumnWidth = 19.6666666667
